$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new package entry as row 6, mirroring the layout of the existing
# rows (Package Name, Package ID, Artifact ID, Version, Type, Date Uploaded).
$ws.Range("A6").Value = "tespack"
$ws.Range("B6").Value = "tespack"
$ws.Range("C6").Value = "wedsf"
$ws.Range("D6").Value = "1.0.0"
$ws.Range("E6").Value = "IFlow"

# "Date Uploaded" is stored as plain text (e.g. "2025-08-08" above is text,
# not a real date), so force the new cell to text as well before assigning
# it - otherwise Excel auto-converts the ISO-formatted string into a date
# serial number. Re-applying the "Normal" style afterwards keeps the cell's
# formatting the same as all of its neighbours.
$dateCell = $ws.Range("F6")
$dateCell.NumberFormat = "@"
$dateCell.Value = "2026-02-04"
$dateCell.Style = "Normal"
